$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 897, shifting existing row 897 (and below) down.
$ws.Rows.Item(897).Insert()

# Populate the newly inserted row 897 with the new data point.
# Force column A to remain plain text (matches the other date-as-text cells)
# instead of being auto-converted into a date serial number.
$ws.Cells.Item(897, 1).NumberFormat = "@"
$ws.Cells.Item(897, 1).Value = "2026/03/01"
$ws.Cells.Item(897, 2).Value = "日"
$ws.Cells.Item(897, 3).Value = 4
$ws.Cells.Item(897, 4).Value = 201
